# AfDD_2023_Annex_Table_Tab21.xlsx - data refresh on "Tab21"
# - Fix mojibake in the Regional Economic Communities legend (A103)
# - Update a handful of recalculated indicator values (F68, C70, and the
#   full "Africa, Fragile States" / "ROW, Fragile States" rows 97-98)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab21")

# --- Fix corrupted accented characters in the footnote text ---
$ws.Range("A103").Value = 'Regional Economic Communities:CEN-SAD = "Community of Sahel-Saharan States";COMESA = "Common Market for Eastern and Southern Africa";EAC = "East African Community";ECCAS = "Economic Community of Central African States";ECOWAS = "Economic Community of West African States";IGAD = "Intergovernmental Authority on Development";SADC = "Southern African Development Community";UMA = "Arab Maghreb Union";PALOP = "Países Africanos de Língua Oficial Portuguesa";ASEAN = "Association of Southeast Asian Nations";MERCOSUR = "Mercado Común del Sur".EU27 = "European Union (27 members)".OECD = "Organisation for Economic Co-operation and Development".'

# --- Minor recalculated value tweaks ---
$ws.Range("F68").Value = 50.126081540933598
$ws.Range("C70").Value = 3.2505543558479202

# --- Row 97: "Africa, Fragile States" ---
$ws.Range("C97").Value = 3.7105815364924402
$ws.Range("D97").Value = 0.097212359640300006
$ws.Range("E97").Value = 11.5621076101411
$ws.Range("F97").Value = 58.328198835185901
$ws.Range("G97").Value = 1.2448215515647401
$ws.Range("H97").Value = 104.63437671918901
$ws.Range("I97").Value = 60329.368784655402
$ws.Range("J97").Value = 215636.77022395501
$ws.Range("K97").Value = 1150102.7246153201
$ws.Range("L97").Value = 23569.139922166702
$ws.Range("M97").Value = 2046984.07795194

# --- Row 98: "ROW, Fragile States" ---
$ws.Range("C98").Value = 3.6181528444481499
$ws.Range("D98").Value = 0.32774794564791998
$ws.Range("E98").Value = 13.009637607434
$ws.Range("F98").Value = 113.85409012802199
$ws.Range("G98").Value = 2.8816169937430098
$ws.Range("H98").Value = 178.37314357998
$ws.Range("I98").Value = 46354.585486524898
$ws.Range("J98").Value = 179716.136546372
$ws.Range("K98").Value = 663537.456072442
$ws.Range("L98").Value = 50321.909004216002
$ws.Range("M98").Value = 1790997.8755076099
